$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row right above the current row 45, shifting all
# existing rows 45-174 down to 46-175 (so the table grows from 174 to 175
# used rows, dimension A1:R174 -> A1:R175).
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row 45 with the new record.
$ws.Range("A45").Value = 9
$ws.Range("B45").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C45").Value = "Metropolitana"
$ws.Range("D45").Value = 44497
$ws.Range("E45").Value = 13
$ws.Range("F45").Value = 300000001
$ws.Range("G45").Value = "Rabanito"
$ws.Range("H45").Value = "Sin especificar"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 5000
$ws.Range("K45").Value = 3500
$ws.Range("L45").Value = 4000
$ws.Range("M45").Value = 3800
$ws.Range("N45").Value = "$/cien unidades (volumen en unidades)"
$ws.Range("O45").Value = "Región Metropolitana"
$ws.Range("P45").Value = 38
$ws.Range("Q45").Value = 100
$ws.Range("R45").Value = "Hortaliza"
